$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new leading row for 2022-Q4 and renumber the
#    0-based index column for the rows that shift down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 36
$summary.Range("D2").Value = 9.41

# Re-use the existing header format (bold / centered / boxed) for the new
# index cell instead of inventing a new style.
$summary.Range("B1").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("A2").Value = 0
$wb.Application.CutCopyMode = $false

for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Add a brand-new "2022-Q4" sheet right after "总计", holding the fund
#    holdings detail for the quarter (mirrors the layout of the other
#    quarterly sheets).
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Header formatting - copy from the same reusable boxed/bold style.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$fundData = @(
    @('013220', '中欧新兴价值一年持有混合A', '36.58', '93.35', '5.83', '2.1326', '5'),
    @('002121', '广发沪港深新起点股票A', '26.30', '88.97', '5.56', '1.4623', '3'),
    @('013221', '中欧新兴价值一年持有混合C', '16.10', '93.35', '5.83', '0.9386', '5'),
    @('014404', '中欧多元价值三年持有混合A', '14.47', '91.65', '5.65', '0.8176', '6'),
    @('010761', '华商甄选回报混合A', '40.49', '70.56', '1.28', '0.5183', '8'),
    @('011856', '安信均衡成长18个月持有混合A', '4.90', '92.55', '9.15', '0.4484', '3'),
    @('005583', '易方达港股通红利灵活配置混合', '7.47', '88.74', '4.65', '0.3474', '4'),
    @('202801', '南方全球精选配置（QDII-FOF）', '17.02', '32.64', '1.67', '0.2842', '5'),
    @('016049', '华商甄选回报混合C', '20.87', '70.56', '1.28', '0.2671', '8'),
    @('008891', '安信价值成长混合A', '2.91', '88.51', '7.57', '0.2203', '3'),
    @('003304', '前海开源沪港深核心资源灵活配置混合A', '3.30', '90.48', '6.50', '0.2145', '9'),
    @('003305', '前海开源沪港深核心资源灵活配置混合C', '3.17', '90.48', '6.50', '0.2060', '9'),
    @('009880', '安信成长动力一年持有期混合', '1.94', '93.22', '8.77', '0.1701', '3'),
    @('008892', '安信价值成长混合C', '2.21', '88.51', '7.57', '0.1673', '3'),
    @('008488', '华商恒益稳健混合', '4.22', '49.85', '2.70', '0.1139', '2'),
    @('011583', '大成港股精选混合（QDII）A', '2.76', '88.26', '3.92', '0.1082', '7'),
    @('007109', '南方沪港深核心优势混合', '1.99', '89.23', '4.75', '0.0945', '4'),
    @('004497', '前海开源多元策略灵活配置混合C', '1.79', '79.66', '4.98', '0.0891', '2'),
    @('012924', '华夏新时代灵活配置混合（QDII）美元现汇', '2.09', '77.57', '4.20', '0.0878', '2'),
    @('012925', '华夏新时代灵活配置混合（QDII）美元现钞', '2.09', '77.57', '4.20', '0.0878', '2'),
    @('012252', '安信宏盈18个月持有混合', '9.14', '25.05', '0.93', '0.0850', '7'),
    @('014405', '中欧多元价值三年持有混合C', '1.36', '91.65', '5.65', '0.0768', '6'),
    @('014746', '贝莱德港股通远景视野混合A', '2.26', '82.52', '3.28', '0.0741', '9'),
    @('004496', '前海开源多元策略灵活配置混合A', '1.30', '79.66', '4.98', '0.0647', '2'),
    @('014747', '贝莱德港股通远景视野混合C', '1.66', '82.52', '3.28', '0.0544', '9'),
    @('014621', '安信楚盈一年持有混合A', '5.80', '22.23', '0.66', '0.0383', '10'),
    @('003243', '上投摩根中国世纪灵活配置混合人民币份额（QDII）', '1.24', '85.53', '3.00', '0.0372', '10'),
    @('003244', '上投摩根中国世纪灵活配置混合美元现钞（QDII）', '1.24', '85.53', '3.00', '0.0372', '10'),
    @('003245', '上投摩根中国世纪灵活配置混合美元现汇（QDII）', '1.24', '85.53', '3.00', '0.0372', '10'),
    @('010024', '广发沪港深新起点股票C', '0.59', '88.97', '5.56', '0.0328', '3'),
    @('004098', '前海开源港股通股息率50强股票', '0.53', '90.79', '5.51', '0.0292', '2'),
    @('011857', '安信均衡成长18个月持有混合C', '0.26', '92.55', '9.15', '0.0238', '3'),
    @('011584', '大成港股精选混合（QDII）C', '0.44', '88.26', '3.92', '0.0172', '7'),
    @('014622', '安信楚盈一年持有混合C', '2.11', '22.23', '0.66', '0.0139', '10'),
    @('501303', '广发恒生中型股指数（LOF）A', '0.24', '90.77', '2.12', '0.0051', '2'),
    @('004996', '广发恒生中型股指数（LOF）C', '0.12', '90.77', '2.12', '0.0025', '2'),
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]

    $q4.Cells.Item($row, 1).Value = $i

    # Columns B-G are stored as plain text (fund code keeps its leading
    # zero, the numeric-looking figures keep their original formatting) -
    # force text entry, then drop the residual "@" number format so the
    # cell is left with the sheet's default (un-styled) look.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q4.Cells.Item($row, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$c - 2]
        $cell.ClearFormats()
    }

    $q4.Cells.Item($row, 8).Value = [int]$rec[6]
}

# Index column (A) uses the same boxed/bold style as the header row and the
# "总计" index column - applied last so it is not wiped by ClearFormats above.
$summary.Range("B1").Copy()
$q4.Range("A2:A37").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false
for ($i = 0; $i -lt $fundData.Length; $i++) {
    $q4.Cells.Item($i + 2, 1).Value = $i
}

# Restore the original active sheet / selection ("总计", cell A1).
$summary.Activate() | Out-Null
$summary.Range("A1").Select() | Out-Null

Write-Output "edit complete"
